$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: B8 and E8 become static values (t-test inputs replaced),
# C8:D8 keep a shared formula (=C5-1 style)
$ws.Range("C8:D8").Formula = "=C5-1"
$ws.Range("B8").Value = 19.489999999999998
$ws.Range("E8").Value = 12.91

# Row 9-10: B9:D10 keep a shared formula (=B6-1 style),
# E9 becomes its own individual formula (=E6-2),
# E10 becomes a static value
$ws.Range("B9:D10").Formula = "=B6-1"
$ws.Range("E9").Formula = "=E6-2"
$ws.Range("E10").Value = 12.96

# Update the active selection to F12
$ws.Range("F12").Select()
